# NV-16 Kha Nhu Huynh 8-2024.xlsx edit
# Commit: "thay doi chien luoc chay multi process. Sua lai template bao cao tong hop co so"
#
# Summary of change on sheet "Đơn phụ phẫu 1":
#  - A new service record (prefix HD-LUXURY, code 618, customer "đường thị út",
#    service "Cắt mí", amount 50000) is inserted as row 2, pushing the
#    existing record (code 619, "mai hồng nương", "Thu cánh mũi", 100000) to
#    row 3, and the "Tổng" row moves to row 4 with updated totals (count 2,
#    amount 150000).
# On sheet "Lương" the dependent totals for "Công phụ phẫu 1 tại SÓC TRĂNG"
# and the rollups that include it are updated from 100000 to 150000.

$wb = $excel.ActiveWorkbook

$wsPhuPhau = $wb.Worksheets.Item("Đơn phụ phẫu 1")

# Insert a new row above the current row 2 data row, so the existing data
# row moves from row 2 -> row 3 and the "Tổng" row moves from row 3 -> row 4.
$wsPhuPhau.Rows.Item(2).Insert()

# Row 2: new service record.
# The "Ngày thực hiện" column holds a literal text date (dd-mm-yyyy), not a
# real date value. Force the cell to Text format before assigning so it
# isn't auto-coerced into a date serial number, then clear the (now
# redundant) explicit formatting so no stray number-format style is left
# behind on the cell.
$wsPhuPhau.Cells.Item(2, 1).Value = "HD-LUXURY"
$wsPhuPhau.Cells.Item(2, 2).Value = 618
$wsPhuPhau.Cells.Item(2, 3).NumberFormat = "@"
$wsPhuPhau.Cells.Item(2, 3).Value = "08-02-2024"
$wsPhuPhau.Cells.Item(2, 3).ClearFormats()
$wsPhuPhau.Cells.Item(2, 4).Value = "SÓC TRĂNG"
$wsPhuPhau.Cells.Item(2, 5).Value = "đường thị út"
$wsPhuPhau.Cells.Item(2, 6).Value = "Cá nhân"
$wsPhuPhau.Cells.Item(2, 7).Value = "Cắt mí"
$wsPhuPhau.Cells.Item(2, 8).Value = "Kha Như Huỳnh "
$wsPhuPhau.Cells.Item(2, 9).Value = 50000

# Row 3: original service record (unchanged values), now shifted down.
$wsPhuPhau.Cells.Item(3, 1).Value = "HD-LUXURY"
$wsPhuPhau.Cells.Item(3, 2).Value = 619
$wsPhuPhau.Cells.Item(3, 3).NumberFormat = "@"
$wsPhuPhau.Cells.Item(3, 3).Value = "08-02-2024"
$wsPhuPhau.Cells.Item(3, 3).ClearFormats()
$wsPhuPhau.Cells.Item(3, 4).Value = "SÓC TRĂNG"
$wsPhuPhau.Cells.Item(3, 5).Value = "mai hồng nương"
$wsPhuPhau.Cells.Item(3, 6).Value = "Cá nhân"
$wsPhuPhau.Cells.Item(3, 7).Value = "Thu cánh mũi"
$wsPhuPhau.Cells.Item(3, 8).Value = "Kha Như Huỳnh "
$wsPhuPhau.Cells.Item(3, 9).Value = 100000

# Row 4: "Tổng" row, now reflecting 2 records / 150000 total.
$wsPhuPhau.Cells.Item(4, 1).Value = "Tổng"
$wsPhuPhau.Cells.Item(4, 2).Value = 2
$wsPhuPhau.Cells.Item(4, 9).Value = 150000

# Update the dependent "Lương" (salary) summary sheet: the "Công phụ phẫu 1
# tại SÓC TRĂNG" line and its downstream rollups move from 100000 to 150000.
$wsLuong = $wb.Worksheets.Item("Lương")
$wsLuong.Cells.Item(29, 2).Value = 150000
$wsLuong.Cells.Item(34, 2).Value = 150000
$wsLuong.Cells.Item(35, 2).Value = 150000
